$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Friday hours for the week of row 13 (F13: 1 -> 2)
$ws.Range("F13").Value = 2

# Add Saturday hours for the week of row 13 (G13: new value 5)
$ws.Range("G13").Value = 5

# Update the selection to match the final saved state (K18:L18, active cell K18)
$ws.Range("K18:L18").Select()

$wb.Save()
